$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: fill in the missing "Meeting" (column A) for the existing British
# Society of Immunology / Belfast poster entry. Copy the formatting (italic
# -free Arial style) from the other "British Society of Immunology" row
# (A28) so we land on the same cell style rather than minting a new one.
$ws.Range("A29").Value = "British Society of Immunology"
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# Row 30: brand-new entry - UK Kidney Week (Invited talk), Oral presentation,
# 2024, Edinburgh. Write column D (Edinburgh) before column A (UK Kidney
# Week) so the new shared strings get appended in that order.
$ws.Range("D30").Value = "Edinburgh"
$ws.Range("B30").Value = "Oral presentation"
$ws.Range("C30").Value = 2024
$ws.Range("A30").Value = "UK Kidney Week (Invited talk)"
$ws.Range("A28").Copy()
$ws.Range("A30").PasteSpecial(-4122)

# Clear clipboard marching-ants/selection state, then move the cursor to
# match the author's last position.
$excel.CutCopyMode = $false
$ws.Range("A32").Select()
